$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 6).Value = 2.18
$ws.Cells.Item(2, 7).Value = 2.26
$ws.Cells.Item(2, 8).Value = 3.2
$ws.Cells.Item(2, 9).Value = 3.4
$ws.Cells.Item(2, 11).Value = 4.2
$ws.Cells.Item(2, 12).Value = 1.3
$ws.Cells.Item(2, 13).Value = 1.04
$ws.Cells.Item(2, 15).Value = 1.2
$ws.Cells.Item(2, 16).Value = 2.54
$ws.Cells.Item(2, 17).Value = 1.6
$ws.Cells.Item(2, 18).Value = 1.62
$ws.Cells.Item(2, 19).Value = 2.5
$ws.Cells.Item(2, 20).Value = 1.57
$ws.Cells.Item(2, 21).Value = 2.58
$ws.Cells.Item(2, 22).Value = 1.42
$ws.Cells.Item(2, 23).Value = 1.79
$ws.Cells.Item(2, 25).Value = 19
$ws.Cells.Item(2, 27).Value = 300
$ws.Cells.Item(2, 28).Value = 15
$ws.Cells.Item(2, 29).Value = 9.8
$ws.Cells.Item(2, 30).Value = 14.5
$ws.Cells.Item(2, 31).Value = 32
$ws.Cells.Item(2, 32).Value = 17
$ws.Cells.Item(2, 34).Value = 14.5
$ws.Cells.Item(2, 35).Value = 36
$ws.Cells.Item(2, 36).Value = 29
$ws.Cells.Item(2, 39).Value = 60
$ws.Cells.Item(2, 40).Value = 11
$ws.Cells.Item(2, 41).Value = 22

# Row 3
$ws.Cells.Item(3, 9).Value = 3.85
$ws.Cells.Item(3, 10).Value = 3.9
$ws.Cells.Item(3, 12).Value = 1.32
$ws.Cells.Item(3, 14).Value = 5.1
$ws.Cells.Item(3, 16).Value = 2.4
$ws.Cells.Item(3, 17).Value = 1.64
$ws.Cells.Item(3, 18).Value = 1.55
$ws.Cells.Item(3, 19).Value = 2.62

# Row 4
$ws.Cells.Item(4, 6).Value = 5.7
$ws.Cells.Item(4, 7).Value = 6.2
$ws.Cells.Item(4, 8).Value = 1.67
$ws.Cells.Item(4, 9).Value = 1.72
$ws.Cells.Item(4, 11).Value = 4.4
$ws.Cells.Item(4, 13).Value = 1.06
$ws.Cells.Item(4, 15).Value = 1.28
$ws.Cells.Item(4, 16).Value = 2.08
$ws.Cells.Item(4, 17).Value = 1.83
$ws.Cells.Item(4, 18).Value = 1.42
$ws.Cells.Item(4, 20).Value = 1.8
$ws.Cells.Item(4, 21).Value = 2.1
$ws.Cells.Item(4, 22).Value = 2.36
$ws.Cells.Item(4, 25).Value = 9.8
$ws.Cells.Item(4, 26).Value = 11
$ws.Cells.Item(4, 27).Value = 18
$ws.Cells.Item(4, 28).Value = 48
$ws.Cells.Item(4, 33).Value = 23
$ws.Cells.Item(4, 34).Value = 20
$ws.Cells.Item(4, 35).Value = 95
$ws.Cells.Item(4, 36).Value = 180
$ws.Cells.Item(4, 38).Value = 300

# Row 5
$ws.Cells.Item(5, 6).Value = 1.86
$ws.Cells.Item(5, 7).Value = 2.02
$ws.Cells.Item(5, 8).Value = 4.5
$ws.Cells.Item(5, 9).Value = 5.3
$ws.Cells.Item(5, 10).Value = 3.4
$ws.Cells.Item(5, 11).Value = 3.8
$ws.Cells.Item(5, 12).Value = 1.47
$ws.Cells.Item(5, 13).Value = 1.08
$ws.Cells.Item(5, 14).Value = 3.2
$ws.Cells.Item(5, 16).Value = 1.74
$ws.Cells.Item(5, 17).Value = 2.16
$ws.Cells.Item(5, 18).Value = 1.27
$ws.Cells.Item(5, 19).Value = 2.18
$ws.Cells.Item(5, 20).Value = 1.03
$ws.Cells.Item(5, 21).Value = 1.03
$ws.Cells.Item(5, 23).Value = 1.99
$ws.Cells.Item(5, 24).Value = 1000
$ws.Cells.Item(5, 28).Value = 1000
$ws.Cells.Item(5, 29).Value = 1000

# Row 6
$ws.Cells.Item(6, 6).Value = 3
$ws.Cells.Item(6, 7).Value = 3.1
$ws.Cells.Item(6, 8).Value = 2.8
$ws.Cells.Item(6, 9).Value = 2.9
$ws.Cells.Item(6, 14).Value = 2.66
$ws.Cells.Item(6, 16).Value = 1.53
$ws.Cells.Item(6, 17).Value = 2.74
$ws.Cells.Item(6, 20).Value = 2.18
$ws.Cells.Item(6, 22).Value = 1.53
$ws.Cells.Item(6, 23).Value = 1.48
$ws.Cells.Item(6, 24).Value = 8.4
$ws.Cells.Item(6, 26).Value = 16.5
$ws.Cells.Item(6, 28).Value = 8.6
$ws.Cells.Item(6, 29).Value = 7
$ws.Cells.Item(6, 30).Value = 14
$ws.Cells.Item(6, 31).Value = 44
$ws.Cells.Item(6, 32).Value = 18
$ws.Cells.Item(6, 33).Value = 14.5
$ws.Cells.Item(6, 34).Value = 24
$ws.Cells.Item(6, 35).Value = 75
$ws.Cells.Item(6, 36).Value = 55
$ws.Cells.Item(6, 37).Value = 48
$ws.Cells.Item(6, 38).Value = 75
$ws.Cells.Item(6, 39).Value = 200
$ws.Cells.Item(6, 40).Value = 70

# Row 7
$ws.Cells.Item(7, 6).Value = 1.07
$ws.Cells.Item(7, 7).Value = 1.81
$ws.Cells.Item(7, 8).Value = 5.4
$ws.Cells.Item(7, 9).Value = 16.5
$ws.Cells.Item(7, 10).Value = 3.45
$ws.Cells.Item(7, 17).Value = 1.62
$ws.Cells.Item(7, 18).Value = 1.24
$ws.Cells.Item(7, 19).Value = 2.72
$ws.Cells.Item(7, 22).Value = 1.06
$ws.Cells.Item(7, 23).Value = 2.22

# Row 8
$ws.Cells.Item(8, 7).Value = 1.48
$ws.Cells.Item(8, 8).Value = 8.6
$ws.Cells.Item(8, 9).Value = 12
$ws.Cells.Item(8, 12).Value = 1.43
$ws.Cells.Item(8, 14).Value = 3.4
$ws.Cells.Item(8, 16).Value = 1.83
$ws.Cells.Item(8, 18).Value = 1.31
$ws.Cells.Item(8, 20).Value = 2.22
$ws.Cells.Item(8, 21).Value = 1.65
$ws.Cells.Item(8, 23).Value = 3.05
$ws.Cells.Item(8, 24).Value = 14
$ws.Cells.Item(8, 26).Value = 95
$ws.Cells.Item(8, 27).Value = 490
$ws.Cells.Item(8, 30).Value = 42
$ws.Cells.Item(8, 31).Value = 960
$ws.Cells.Item(8, 33).Value = 13
$ws.Cells.Item(8, 41).Value = 400

# Row 9
$ws.Cells.Item(9, 6).Value = 2.96
$ws.Cells.Item(9, 7).Value = 3.05
$ws.Cells.Item(9, 8).Value = 2.82
$ws.Cells.Item(9, 9).Value = 2.9
$ws.Cells.Item(9, 12).Value = 1.6
$ws.Cells.Item(9, 14).Value = 2.7
$ws.Cells.Item(9, 15).Value = 1.56
$ws.Cells.Item(9, 21).Value = 1.81
$ws.Cells.Item(9, 22).Value = 1.52
$ws.Cells.Item(9, 27).Value = 48
$ws.Cells.Item(9, 41).Value = 50

# Row 10
$ws.Cells.Item(10, 6).Value = 1.37
$ws.Cells.Item(10, 7).Value = 1.39
$ws.Cells.Item(10, 8).Value = 11
$ws.Cells.Item(10, 9).Value = 13
$ws.Cells.Item(10, 11).Value = 5.4
$ws.Cells.Item(10, 12).Value = 1.45
$ws.Cells.Item(10, 16).Value = 1.81
$ws.Cells.Item(10, 17).Value = 2.14
$ws.Cells.Item(10, 20).Value = 2.6
$ws.Cells.Item(10, 23).Value = 3.55
$ws.Cells.Item(10, 26).Value = 120
$ws.Cells.Item(10, 28).Value = 6
$ws.Cells.Item(10, 29).Value = 11.5
$ws.Cells.Item(10, 35).Value = 1000
$ws.Cells.Item(10, 39).Value = 450
$ws.Cells.Item(10, 40).Value = 8
$ws.Cells.Item(10, 41).Value = 680
